$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00001292064567892659
$ws.Range("C2").Value = 1.655778082260271
$ws.Range("D2").Value = 57107556.33100624
$ws.Range("E2").Value = 85231193291209616
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 85231193348317168
